# Updates cryptos list values (price & volume) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.388.40'
$ws.Range("E2").Value = '  -7.38%  '
$ws.Range("D3").Value = '3.308.71'
$ws.Range("E3").Value = '  -4.88%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '''559.18'
$ws.Range("E5").Value = '  -4.45%  '
$ws.Range("D6").Value = '''128.10'
$ws.Range("E6").Value = '  -2.48%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.313.08'
$ws.Range("E8").Value = '  -4.73%  '
$ws.Range("D9").Value = '''0.467'
$ws.Range("E9").Value = '  -3.10%  '
$ws.Range("D10").Value = '''7.39'
$ws.Range("E10").Value = '  -4.23%  '
$ws.Range("E11").Value = '  -5.91%  '
$ws.Range("D12").Value = '''0.370'
$ws.Range("E12").Value = '  -4.20%  '
$ws.Range("D13").Value = '3.863.01'
$ws.Range("E13").Value = '  -5.12%  '
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = '3.287.23'
$ws.Range("E15").Value = '  -5.44%  '
$ws.Range("E16").Value = '  -6.74%  '
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = '''23.94'
$ws.Range("E17").Value = '  -4.80%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '59.581.47'
$ws.Range("E18").Value = '  -7.07%  '
$ws.Range("D19").Value = '''5.63'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").Value = '''13.22'
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").Value = '''8.89'
$ws.Range("E21").Value = '  -10.89%  '
$ws.Range("D22").Value = '''350.30'
$ws.Range("E22").Value = '  -9.17%  '
$ws.Range("D23").Value = '''0.552'
$ws.Range("E23").Value = '  -2.78%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").Value = '3.432.11'
$ws.Range("E25").Value = '  -5.15%  '
$ws.Range("D26").Value = '''68.33'
$ws.Range("E26").Value = '  -8.28%  '
$ws.Range("D27").Value = '''0.0000108'
$ws.Range("E27").Value = '  -3.63%  '
$ws.Range("D28").Value = '''0.998'
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = '''7.25'
$ws.Range("E29").Value = '  +2.30%  '
$ws.Range("D30").Value = '''1.42'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  -2.48%  '
$ws.Range("D32").Value = '''0.150'
$ws.Range("E32").Value = '  -3.73%  '
$ws.Range("E33").Value = '  -6.23%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("D35").Value = '3.333.33'
$ws.Range("E35").Value = '  -4.97%  '
$ws.Range("D36").Value = '''22.73'
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("D37").Value = '''5.26'
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("D38").Value = '''6.76'
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("D39").Value = '''1.47'
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("D40").Value = '''157.08'
$ws.Range("E40").Value = '  -3.84%  '
$ws.Range("D41").Value = '''0.0745'
$ws.Range("E41").Value = '  -3.91%  '
$ws.Range("D42").Value = '''0.998'
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").Value = '''40.39'
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("D44").Value = '''0.741'
$ws.Range("E44").Value = '  -7.05%  '
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("E46").Value = '  +2.54%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''22.54'
$ws.Range("E47").Value = '  -4.37%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '''1.53'
$ws.Range("E48").Value = '  -5.34%  '
$ws.Range("D49").Value = '''6.70'
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("D50").Value = '''21.66'
$ws.Range("E50").Value = '  +5.71%  '
$ws.Range("B51").Value = 'LidoDAOToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D51").Value = '''2.32'
$ws.Range("E51").Value = '  +10.99%  '

# Cells that needed a leading quote to stay text (since their new value looks like a
# plain number) get their style reset to Normal so no extra formatting/style is left behind.
foreach ($addr in @("D4", "D5", "D6", "D9", "D10", "D12", "D17", "D19", "D20", "D21", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D32", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D50", "D51")) {
    $ws.Range($addr).Style = "Normal"
}
